$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1. Refresh the cached "datetimeFigureOut" footer field text on the
#    slide master and every slide layout (11-7-2016 -> 15-11-2016).
# ------------------------------------------------------------------
$newDate = "15-11-2016"

$master = $p.Slides.Item(1).Master

# Update the date placeholder on the slide master itself.
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq "11-7-2016") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Update the date placeholder on every custom (slide) layout.
$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "11-7-2016") {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ------------------------------------------------------------------
# 2. Slide 1: "Install Python 3.5 from www.python.org" becomes
#    "Install Python 3.5 or 3.6 from www.python.org" (support for
#    yield from - Python 3.6 is now also recommended).
# ------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -like "Install Python 3.5 from www.python.org*") {
            # "from www.python.org" starts right after "Install Python 3.5 " (19 chars).
            # Grab the "from " word (5 chars) and insert "or 3.6 " just before it so the
            # inserted text merges with "from " into a single run, matching:
            #   "Install Python 3.5 " | "or 3.6 from " | "www.python.org"
            $fromRange = $tr.Characters(20, 5)
            [void]$fromRange.InsertBefore("or 3.6 ")
        }
    }
}
